$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.323.31'
$ws.Range('E2').Value = '  +2.56%  '
$ws.Range('D3').Value = '2.423.63'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.07'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +2.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.72'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.95%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +1.96%  '
$ws.Range('D9').Value = '2.424.37'
$ws.Range('E9').Value = '  +3.31%  '
$ws.Range('E10').Value = '  +4.60%  '
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('E12').Value = '  +1.65%  '
$ws.Range('E13').Value = '  +2.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.26'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +6.29%  '
$ws.Range('E15').Value = '  +9.15%  '
$ws.Range('D16').Value = '2.861.34'
$ws.Range('E16').Value = '  +3.22%  '
$ws.Range('D17').Value = '62.160.59'
$ws.Range('E17').Value = '  +2.30%  '
$ws.Range('D18').Value = '2.421.94'
$ws.Range('E18').Value = '  +3.14%  '
$ws.Range('E19').Value = '  +4.19%  '
$ws.Range('E20').Value = '  +1.87%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.57'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.76%  '
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('E24').Value = '  +6.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '64.82'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.41%  '
$ws.Range('E26').Value = '  +9.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '572.61'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +15.12%  '
$ws.Range('D28').Value = '2.538.61'
$ws.Range('E28').Value = '  +3.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.40'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +5.50%  '
$ws.Range('E31').Value = '  +9.55%  '
$ws.Range('E32').Value = '  +6.05%  '
$ws.Range('E33').Value = '  +2.00%  '
$ws.Range('E34').Value = '  +3.92%  '
$ws.Range('E35').Value = '  +4.91%  '
$ws.Range('E36').Value = '  +9.20%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  +5.16%  '
$ws.Range('E39').Value = '  +2.65%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.88'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.58%  '
$ws.Range('B41').Value = 'EthereumClassic'
$ws.Range('C41').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '18.77'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '150.02'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.00%  '
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '41.66'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +2.56%  '
$ws.Range('E45').Value = '  +14.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '150.96'
$ws.Range('D46').ClearFormats()
$ws.Range('E47').Value = '  +2.33%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0541'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.77%  '
$ws.Range('E49').Value = '  +6.94%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.589'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +3.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0229'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +3.69%  '
